$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: fill in Gross Profit formulas for I6:L6 (J6:L6 as a shared formula group like H6 pattern)
$ws.Range("I6").Formula = "=I4-I5"
$ws.Range("J6:L6").Formula = "=J4-J5"

# Apply matching currency-style formatting (bold, top border, centered) like the rest of row 6
$ws.Range("I6:L6").Font.Bold = $true
$ws.Range("I6:L6").Borders.Item(8).LineStyle = 1
$ws.Range("I6:L6").NumberFormat = "[$$-409]#,##0.00"
$ws.Range("I6:L6").HorizontalAlignment = -4108

# Row 8: replace the "Put Formula Here" placeholder with the actual average formula
$ws.Range("C8").Formula = "=AVERAGE(C6:L6)"
$ws.Range("C8").Font.Bold = $false
$ws.Range("C8").NumberFormat = "_(""$""* #,##0.00_);_(""$""* \(#,##0.00\);_(""$""* ""-""??_);_(@_)"

# Match the final selection left behind in the saved workbook
$ws.Range("C8").Select() | Out-Null
